# Quarterly database update: roll the 10-quarter reporting window forward by
# one quarter (drop "Q2 ended 1399/06", add new "Q4 ended 1401/12") and
# refresh the read_price-derived figures for the shifted window.
#
# Layout: row 8 and row 24 hold the ten quarter-end labels (columns E:N).
# Rows 16/17/19/20 (cost table) and rows 26/27 (headcount table) hold the
# corresponding ten data points for the same column window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

$dataCols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# Header rows: both quarter-label strips shift one column to the left, with
# the brand-new quarter label landing in column N.
$headerRows = @(8, 24)
foreach ($row in $headerRows) {
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range($dataCols[$i] + $row).Value = $quarterLabels[$i]
    }
}

# Data rows: each row's ten values shift one column to the left, dropping the
# oldest quarter's figure and appending the newly reported quarter's figure
# in column N.
$rowData = @{
    16 = @(461, 289, 327, 586, 281, 444, 376, 395, 480, 529)
    17 = @(11235, 24145, 35145, 25881, 35796, 31870, 43786, 42594, 53832, 88833)
    19 = @(20861, 22824, 12069, 27358, 10818, 124615, 16294, 38438, 16571, 37426)
    20 = @(32557, 47258, 47541, 53825, 46895, 156929, 60456, 81427, 70883, 126788)
    26 = @(269, 271, 272, 273, 272, 273, 264, 270, 270, 272)
    27 = @(265, 258, 255, 255, 251, 253, 252, 247, 244, 243)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range($dataCols[$i] + $row).Value = $values[$i]
    }
}
